# Update Pais sheet: refresh COVID stats + re-rank a few countries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 01:52"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 364833
$ws.Range("C4").Value = 28160
$ws.Range("D4").Value = 19536
$ws.Range("E4").Value = 334438
$ws.Range("F4").Value = 8878
$ws.Range("G4").Value = 1243
$ws.Range("H4").Value = 10859

# Row 18: 'Brasil' -> 'Brasil'
$ws.Range("B18").Value = 12183
$ws.Range("C18").Value = 929
$ws.Range("D18").Value = 127
$ws.Range("E18").Value = 11492
$ws.Range("F18").Value = 296
$ws.Range("G18").Value = 78
$ws.Range("H18").Value = 564

# Row 46: 'Emiratos Arabes Unidos' -> 'Panama'
$ws.Range("A46").Value = "Panama"
$ws.Range("B46").Value = 2100
$ws.Range("C46").Value = 299
$ws.Range("D46").Value = 14
$ws.Range("E46").Value = 2031
$ws.Range("F46").Value = 88
$ws.Range("G46").Value = 9
$ws.Range("H46").Value = 55

# Row 47: 'Panama' -> 'Emiratos Arabes Unidos'
$ws.Range("A47").Value = "Emiratos Arabes Unidos"
$ws.Range("B47").Value = 2076
$ws.Range("C47").Value = 277
$ws.Range("D47").Value = 167
$ws.Range("E47").Value = 1898
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 11

# Row 88: 'Uruguay' -> 'Uruguay'
$ws.Range("B88").Value = 415
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 123
$ws.Range("E88").Value = 286
$ws.Range("F88").Value = 14
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 6

# Row 99: 'San Marino' -> 'San Marino'
$ws.Range("B99").Value = 277
$ws.Range("C99").Value = 11
$ws.Range("D99").Value = 35
$ws.Range("E99").Value = 210
$ws.Range("F99").Value = 14
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 32

# Row 115: 'Consejo Danes para los Refugiados' -> 'Mayotte'
$ws.Range("A115").Value = "Mayotte"
$ws.Range("B115").Value = 164
$ws.Range("C115").Value = 17
$ws.Range("D115").Value = 15
$ws.Range("E115").Value = 147
$ws.Range("F115").Value = 3
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 2

# Row 116: 'Kenia' -> 'Consejo Danes para los Refugiados'
$ws.Range("A116").Value = "Consejo Danes para los Refugiados"
$ws.Range("B116").Value = 161
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 5
$ws.Range("E116").Value = 138
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 18

# Row 117: 'Martinica' -> 'Kenia'
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 158
$ws.Range("C117").Value = 16
$ws.Range("D117").Value = 4
$ws.Range("E117").Value = 148
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 6

# Row 118: 'Mayotte' -> 'Martinica'
$ws.Range("A118").Value = "Martinica"
$ws.Range("B118").Value = 151
$ws.Range("C118").Value = 2
$ws.Range("D118").Value = 50
$ws.Range("E118").Value = 97
$ws.Range("F118").Value = 20
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 4

# Row 119: 'Isla de Man' -> 'Guadalupe'
$ws.Range("A119").Value = "Guadalupe"
$ws.Range("B119").Value = 139
$ws.Range("C119").Value = 4
$ws.Range("D119").Value = 31
$ws.Range("E119").Value = 101
$ws.Range("F119").Value = 14
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 7

# Row 120: 'Guadalupe' -> 'Isla de Man'
$ws.Range("A120").Value = "Isla de Man"
$ws.Range("B120").Value = 139
$ws.Range("C120").Value = 12
$ws.Range("D120").Value = 55
$ws.Range("E120").Value = 83
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1

# Row 151: 'Guam' -> 'Bahamas'
$ws.Range("A151").Value = "Bahamas"
$ws.Range("B151").Value = 33
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 5
$ws.Range("E151").Value = 23
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 5

# Row 152: 'San Martin (Parte Francesa)' -> 'Guam'
$ws.Range("A152").Value = "Guam"
$ws.Range("B152").Value = 32
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 0
$ws.Range("E152").Value = 31
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1

# Row 153: 'Eritrea' -> 'San Martin (Parte Francesa)'
$ws.Range("A153").Value = "San Martin (Parte Francesa)"
$ws.Range("B153").Value = 32
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 7
$ws.Range("E153").Value = 23
$ws.Range("F153").Value = 6
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 2

# Row 154: 'Guyana' -> 'Eritrea'
$ws.Range("A154").Value = "Eritrea"
$ws.Range("B154").Value = 31
$ws.Range("C154").Value = 2
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 31
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 0

# Row 155: 'Bahamas' -> 'Guyana'
$ws.Range("A155").Value = "Guyana"
$ws.Range("B155").Value = 31
$ws.Range("C155").Value = 7
$ws.Range("D155").Value = 8
$ws.Range("E155").Value = 19
$ws.Range("F155").Value = 8
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 4
